$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) "json" list (column M) gains a new, alphabetically-sorted entry:
#    storeKeys(json,jsonpath,var) -- inserted right before storeValue(...)
$ws.Range("M16").Insert(-4121)
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# 2) "target" list (column A) loses its "text" entry (row 25), shifting
#    web/webalert/webcookie/ws/ws.async/xml up by one row.
$ws.Range("A25").Delete(-4162)

# 3) The block of columns Z:AE (web, webalert, webcookie, ws, ws.async, xml)
#    shifts one column to the left (into Y:AD), effectively deleting the old
#    stand-alone column Y (which only backed the single-cell "text" range).
$ws.Range("Y:Y").Delete()

# 4) Update the defined names so their ranges reflect the edits above.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"
